$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Permits Filed for 25-21 46th Street in Astoria, Queens"
$ws.Range("B2").Value = "https://newyorkyimby.com/2025/10/permits-filed-for-25-21-46th-street-in-astoria-queens.html"
$ws.Range("C2").Value = 'Permits have been filed for a three-story residential building at 25-21 46th Street in <a href="https://newyorkyimby.com/neighborhoods/astoria">Astoria</a>, Queens. Located between 25th and 28th Avenues, the lot is near the 46th Street subway station, served by the E, F, M, and R trains. Vincent Maimone of Artistic Design Corp. is listed as the owner behind the applications.'
$ws.Range("D2").Value = "2025-10-17T10:30:44+00:00"
$ws.Range("E2").Value = "Fri, 17 Oct 2025 10:30:44 +0000"
$ws.Range("F2").Value = "YIMBY"
$ws.Range("G2").Value = "YIMBY - Astoria"
$ws.Range("H2").Value = ""
